$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.288.79"
$ws.Range("E2").Value = "'  +3.12%  "
$ws.Range("D3").Value = "'1.816.22"
$ws.Range("E3").Value = "'  +4.14%  "
$ws.Range("D5").Value = "'328.29"
$ws.Range("E5").Value = "'  +2.01%  "
$ws.Range("E6").Value = "'  +0.03%  "
$ws.Range("E7").Value = "'  +3.10%  "
$ws.Range("D8").Value = "'0.3669"
$ws.Range("E8").Value = "'  +2.37%  "
$ws.Range("D9").Value = "'44.90"
$ws.Range("E9").Value = "'  -1.31%  "
$ws.Range("D10").Value = "'0.07687"
$ws.Range("E10").Value = "'  +3.44%  "
$ws.Range("D11").Value = "'1.142"
$ws.Range("E11").Value = "'  +2.62%  "
$ws.Range("E12").Value = "'  +0.08%  "
$ws.Range("D13").Value = "'22.17"
$ws.Range("E13").Value = "'  +3.36%  "
$ws.Range("D14").Value = "'6.309"
$ws.Range("E14").Value = "'  +3.12%  "
$ws.Range("D15").Value = "'7.542"
$ws.Range("E15").Value = "'  +4.86%  "
$ws.Range("D16").Value = "'1.829.95"
$ws.Range("E16").Value = "'  +4.97%  "
$ws.Range("D17").Value = "'93.13"
$ws.Range("E17").Value = "'  +5.85%  "
$ws.Range("D18").Value = "'0.00001082"
$ws.Range("E18").Value = "'  +1.48%  "
$ws.Range("D19").Value = "'0.06534"
$ws.Range("E19").Value = "'  +6.88%  "
$ws.Range("E20").Value = "'  +0.01%  "
$ws.Range("D21").Value = "'17.51"
$ws.Range("E21").Value = "'  +3.82%  "
$ws.Range("D22").Value = "'6.264"
$ws.Range("E22").Value = "'  +2.62%  "
$ws.Range("D23").Value = "'28.315.85"
$ws.Range("E23").Value = "'  +3.10%  "
$ws.Range("D24").Value = "'11.66"
$ws.Range("E24").Value = "'  +1.53%  "
$ws.Range("D25").Value = "'2.049"
$ws.Range("E25").Value = "'  -12.45%  "
$ws.Range("D26").Value = "'162.30"
$ws.Range("E26").Value = "'  +6.38%  "
$ws.Range("D27").Value = "'20.73"
$ws.Range("E27").Value = "'  +1.68%  "
$ws.Range("D28").Value = "'2.033.86"
$ws.Range("E28").Value = "'  +4.83%  "
$ws.Range("D29").Value = "'2.302"
$ws.Range("E29").Value = "'  -3.26%  "
$ws.Range("D30").Value = "'128.73"
$ws.Range("E30").Value = "'  +2.15%  "
$ws.Range("E31").Value = "'  +1.76%  "
$ws.Range("D32").Value = "'5.964"
$ws.Range("E32").Value = "'  +5.11%  "
$ws.Range("D33").Value = "'0.09209"
$ws.Range("E33").Value = "'  +1.03%  "
$ws.Range("D34").Value = "'3.494"
$ws.Range("E34").Value = "'  -3.58%  "
$ws.Range("E35").Value = "'  +2.24%  "
$ws.Range("D36").Value = "'0.02351"
$ws.Range("E36").Value = "'  +2.39%  "
$ws.Range("D37").Value = "'0.2181"
$ws.Range("E37").Value = "'  +2.00%  "
$ws.Range("D38").Value = "'5.199"
$ws.Range("E38").Value = "'  +2.26%  "
$ws.Range("D39").Value = "'0.6580"
$ws.Range("E39").Value = "'  +2.88%  "
$ws.Range("D40").Value = "'0.06209"
$ws.Range("E40").Value = "'  +2.62%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.194"
$ws.Range("E41").Value = "'  +0.39%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.137"
$ws.Range("E42").Value = "'  +3.04%  "
$ws.Range("D43").Value = "'1.436"
$ws.Range("E43").Value = "'  +0.84%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "'  +0.01%  "
$ws.Range("D45").Value = "'13.90"
$ws.Range("E45").Value = "'  +1.79%  "
$ws.Range("D46").Value = "'0.6121"
$ws.Range("E46").Value = "'  +4.19%  "
$ws.Range("D47").Value = "'3.755"
$ws.Range("E47").Value = "'  +1.15%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'125.98"
$ws.Range("E48").Value = "'  +0.66%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.025"
$ws.Range("E49").Value = "'  +4.09%  "
$ws.Range("D50").Value = "'1.159"
$ws.Range("E50").Value = "'  +4.09%  "
$ws.Range("D51").Value = "'0.07010"
$ws.Range("E51").Value = "'  +2.54%  "
